{"js": "// Update the worksheet date and the two-digit multiplication problems.\nconst replacements = [\n    { old: \"2025-11-10 Monday\", new: \"2025-11-11 Tuesday\" },\n    { old: \"99\u00d762=\", new: \"24\u00d780=\" },\n    { old: \"75\u00d797=\", new: \"32\u00d741=\" },\n    { old: \"46\u00d715=\", new: \"72\u00d738=\" },\n    { old: \"36\u00d747=\", new: \"60\u00d724=\" },\n    { old: \"72\u00d759=\", new: \"65\u00d765=\" },\n    { old: \"30\u00d772=\", new: \"32\u00d789=\" },\n    { old: \"13\u00d730=\", new: \"93\u00d763=\" },\n    { old: \"73\u00d753=\", new: \"47\u00d752=\" },\n    { old: \"86\u00d729=\", new: \"45\u00d747=\" },\n    { old: \"92\u00d772=\", new: \"78\u00d741=\" },\n    { old: \"19\u00d726=\", new: \"95\u00d742=\" },\n    { old: \"29\u00d726=\", new: \"52\u00d740=\" },\n    { old: \"17\u00d731=\", new: \"17\u00d752=\" },\n    { old: \"55\u00d734=\", new: \"24\u00d785=\" },\n    { old: \"75\u00d731=\", new: \"86\u00d714=\" },\n    { old: \"68\u00d718=\", new: \"39\u00d767=\" },\n    { old: \"68\u00d756=\", new: \"23\u00d756=\" },\n    { old: \"52\u00d772=\", new: \"48\u00d740=\" },\n    { old: \"57\u00d745=\", new: \"34\u00d723=\" },\n    { old: \"75\u00d716=\", new: \"28\u00d756=\" },\n    { old: \"22\u00d711=\", new: \"26\u00d755=\" },\n    { old: \"59\u00d763=\", new: \"87\u00d729=\" },\n    { old: \"33\u00d749=\", new: \"59\u00d778=\" },\n    { old: \"71\u00d745=\", new: \"53\u00d732=\" },\n    { old: \"21\u00d773=\", new: \"46\u00d723=\" },\n];\n\nconst body = context.document.body;\n\nfor (const pair of replacements) {\n    const results = body.search(pair.old, { matchCase: true });\n    results.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(pair.new, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "# Update the worksheet date and the two-digit multiplication problems.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2025-11-10 Monday\"; New = \"2025-11-11 Tuesday\" },\n    @{ Old = \"99\u00d762=\"; New = \"24\u00d780=\" },\n    @{ Old = \"75\u00d797=\"; New = \"32\u00d741=\" },\n    @{ Old = \"46\u00d715=\"; New = \"72\u00d738=\" },\n    @{ Old = \"36\u00d747=\"; New = \"60\u00d724=\" },\n    @{ Old = \"72\u00d759=\"; New = \"65\u00d765=\" },\n    @{ Old = \"30\u00d772=\"; New = \"32\u00d789=\" },\n    @{ Old = \"13\u00d730=\"; New = \"93\u00d763=\" },\n    @{ Old = \"73\u00d753=\"; New = \"47\u00d752=\" },\n    @{ Old = \"86\u00d729=\"; New = \"45\u00d747=\" },\n    @{ Old = \"92\u00d772=\"; New = \"78\u00d741=\" },\n    @{ Old = \"19\u00d726=\"; New = \"95\u00d742=\" },\n    @{ Old = \"29\u00d726=\"; New = \"52\u00d740=\" },\n    @{ Old = \"17\u00d731=\"; New = \"17\u00d752=\" },\n    @{ Old = \"55\u00d734=\"; New = \"24\u00d785=\" },\n    @{ Old = \"75\u00d731=\"; New = \"86\u00d714=\" },\n    @{ Old = \"68\u00d718=\"; New = \"39\u00d767=\" },\n    @{ Old = \"68\u00d756=\"; New = \"23\u00d756=\" },\n    @{ Old = \"52\u00d772=\"; New = \"48\u00d740=\" },\n    @{ Old = \"57\u00d745=\"; New = \"34\u00d723=\" },\n    @{ Old = \"75\u00d716=\"; New = \"28\u00d756=\" },\n    @{ Old = \"22\u00d711=\"; New = \"26\u00d755=\" },\n    @{ Old = \"59\u00d763=\"; New = \"87\u00d729=\" },\n    @{ Old = \"33\u00d749=\"; New = \"59\u00d778=\" },\n    @{ Old = \"71\u00d745=\"; New = \"53\u00d732=\" },\n    @{ Old = \"21\u00d773=\"; New = \"46\u00d723=\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($pair.Old, $false, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
